# Regenerate merged AHB files: rename header labels from the previous
# "_old" / "_new" diff-suffix convention to the explicit format-version
# suffixes "_FV2410" / "_FV2504", turn the data range into a real Excel
# Table (so AutoFilter + structured references work), and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row labels (row 1, columns A:U) -----------------
# A1:J1  => "<Label>_old"  -> "<Label>_FV2410"
# K1     => "diff"          (unchanged)
# L1:U1  => "<Label>_new"  -> "<Label>_FV2504"
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2
    if ($col -le 10) {
        $cell.Value = $text -replace "_old$", "_FV2410"
    } elseif ($col -ge 12) {
        $cell.Value = $text -replace "_new$", "_FV2504"
    }
}

# --- 2. Convert the used range into a native Excel table -------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$rng = $ws.Range("A1:U" + $lastRow)
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row -----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
